# Append-style scrape refresh: 2025-10-21 06:35 JST
# The oldest listed item ("MT4/MT5用FX自動売買システムの開発者募集", row 6) has
# fallen out of the tracked window; every remaining row moves up by one and
# the capture timestamp in column A is refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the row that is no longer present in this scrape; Excel shifts rows
# 7-9 up into 6-8 and keeps their existing hyperlink relationships attached
# to the (now shifted) F-column cells, just like the source tool does.
$ws.Rows.Item(6).Delete()

# Refresh the "taken at" timestamp for every remaining data row.
$newStamp = "2025-10-21 06:35:04"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newStamp
}

# The row that used to be r9 is gone, but the worksheet's hyperlink list
# still carries a stale entry for it. Rebuild the hyperlink list so it only
# references the rows that still exist (F2:F8), restoring the exact
# relationship targets that were already wired to those cells.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5415908")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5417295")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5408664")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5417433")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5417377")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5371747")
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5417267")
